$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A202").Value = '(Intercept)'
$ws.Range("B202").Value = [double]'4.973137893762494e-116'
$ws.Range("C202").Value = [double]'22.91115374961466'
$ws.Range("D202").Value = [double]'-11.58805980183081'
$ws.Range("E202").Value = [double]'4.737422236229271e-31'
$ws.Range("F202").Value = [double]'1.435448979988836e-135'
$ws.Range("G202").Value = [double]'1.464009066951592e-96'
$ws.Range("H202").Value = 'TZP'

$ws.Range("A203").Value = 'Year'
$ws.Range("B203").Value = [double]'1.140946738498609'
$ws.Range("C203").Value = [double]'0.01135736599391953'
$ws.Range("D203").Value = [double]'11.60994461152569'
$ws.Range("E203").Value = [double]'3.668531608164512e-31'
$ws.Range("F203").Value = [double]'1.115872250963904'
$ws.Range("G203").Value = [double]'1.166679442374917'
$ws.Range("H203").Value = 'TZP'

$ws.Range("A204").Value = 'Specimen_typeother'
$ws.Range("B204").Value = [double]'0.8461702789920515'
$ws.Range("C204").Value = [double]'0.09756391046152789'
$ws.Range("D204").Value = [double]'-1.712053805906441'
$ws.Range("E204").Value = [double]'0.08688675321192382'
$ws.Range("F204").Value = [double]'0.6987173047850542'
$ws.Range("G204").Value = [double]'1.024296154700074'
$ws.Range("H204").Value = 'TZP'

$ws.Range("A205").Value = 'Specimen_typeRespiratory'
$ws.Range("B205").Value = [double]'0.5761503972320655'
$ws.Range("C205").Value = [double]'0.10150406254931'
$ws.Range("D205").Value = [double]'-5.432162341236693'
$ws.Range("E205").Value = [double]'5.567523586273713e-08'
$ws.Range("F205").Value = [double]'0.4719276926138439'
$ws.Range("G205").Value = [double]'0.7026080629646728'
$ws.Range("H205").Value = 'TZP'

$ws.Range("A206").Value = 'Specimen_typeUrine'
$ws.Range("B206").Value = [double]'0.5411927767328973'
$ws.Range("C206").Value = [double]'0.07848801358758777'
$ws.Range("D206").Value = [double]'-7.822592284750599'
$ws.Range("E206").Value = [double]'5.174649281704904e-15'
$ws.Range("F206").Value = [double]'0.4640904728348101'
$ws.Range("G206").Value = [double]'0.6313126397698415'
$ws.Range("H206").Value = 'TZP'

$ws.Range("A207").Value = 'Specimen_typeWound & soft tissues'
$ws.Range("B207").Value = [double]'0.9277326279276334'
$ws.Range("C207").Value = [double]'0.08852634521799822'
$ws.Range("D207").Value = [double]'-0.8473376369715753'
$ws.Range("E207").Value = [double]'0.3968069539696141'
$ws.Range("F207").Value = [double]'0.7799152046028212'
$ws.Range("G207").Value = [double]'1.103513528417497'
$ws.Range("H207").Value = 'TZP'

$ws.Range("A208").Value = 'HospitalCHBH'
$ws.Range("B208").Value = [double]'0.2361244119590601'
$ws.Range("C208").Value = [double]'0.2180553316314181'
$ws.Range("D208").Value = [double]'-6.619404500428351'
$ws.Range("E208").Value = [double]'3.606487429587598e-11'
$ws.Range("F208").Value = [double]'0.1531521194390464'
$ws.Range("G208").Value = [double]'0.3603607512938581'
$ws.Range("H208").Value = 'TZP'

$ws.Range("A209").Value = 'HospitalCNGMO'
$ws.Range("B209").Value = [double]'0.7560454158769172'
$ws.Range("C209").Value = [double]'0.4565262054970577'
$ws.Range("D209").Value = [double]'-0.612569064663486'
$ws.Range("E209").Value = [double]'0.5401613183668608'
$ws.Range("F209").Value = [double]'0.3028431187533294'
$ws.Range("G209").Value = [double]'1.851048722156757'
$ws.Range("H209").Value = 'TZP'

$ws.Range("A210").Value = 'HospitalRabta'
$ws.Range("B210").Value = [double]'1.827064742148086'
$ws.Range("C210").Value = [double]'0.1831295329717927'
$ws.Range("D210").Value = [double]'3.291171572944333'
$ws.Range("E210").Value = [double]'0.0009977104181039741'
$ws.Range("F210").Value = [double]'1.275908334665787'
$ws.Range("G210").Value = [double]'2.617062211604711'
$ws.Range("H210").Value = 'TZP'

$ws.Range("A211").Value = 'HospitalTCB'
$ws.Range("B211").Value = [double]'0.7032866023923675'
$ws.Range("C211").Value = [double]'0.1627257239024766'
$ws.Range("D211").Value = [double]'-2.163092454236802'
$ws.Range("E211").Value = [double]'0.03053406936240863'
$ws.Range("F211").Value = [double]'0.5104589954785741'
$ws.Range("G211").Value = [double]'0.9665934430742568'
$ws.Range("H211").Value = 'TZP'

$ws.Range("A212").Value = 'Ward_ED_ICUED'
$ws.Range("B212").Value = [double]'0.1777510440476676'
$ws.Range("C212").Value = [double]'0.2257028950393169'
$ws.Range("D212").Value = [double]'-7.653297208833935'
$ws.Range("E212").Value = [double]'1.958904520756977e-14'
$ws.Range("F212").Value = [double]'0.1133837245093725'
$ws.Range("G212").Value = [double]'0.274943485114469'
$ws.Range("H212").Value = 'TZP'

$ws.Range("A213").Value = 'Ward_ED_ICUOther'
$ws.Range("B213").Value = [double]'0.2434847467182416'
$ws.Range("C213").Value = [double]'0.1580536060162385'
$ws.Range("D213").Value = [double]'-8.938112934232164'
$ws.Range("E213").Value = [double]'3.958723971324885e-19'
$ws.Range("F213").Value = [double]'0.1783133827609263'
$ws.Range("G213").Value = [double]'0.331537146185214'
$ws.Range("H213").Value = 'TZP'

$ws.Range("A214").Value = 'HospitalCHBH:Ward_ED_ICUED'
$ws.Range("B214").Value = [double]'3.315452634775907'
$ws.Range("C214").Value = [double]'0.4154416661812452'
$ws.Range("D214").Value = [double]'2.885108195200897'
$ws.Range("E214").Value = [double]'0.003912791641065692'
$ws.Range("F214").Value = [double]'1.430144717484845'
$ws.Range("G214").Value = [double]'7.348941456002666'
$ws.Range("H214").Value = 'TZP'

$ws.Range("A215").Value = 'HospitalCNGMO:Ward_ED_ICUED'
$ws.Range("H215").Value = 'TZP'

$ws.Range("A216").Value = 'HospitalRabta:Ward_ED_ICUED'
$ws.Range("H216").Value = 'TZP'

$ws.Range("A217").Value = 'HospitalTCB:Ward_ED_ICUED'
$ws.Range("B217").Value = [double]'1.285038299948908'
$ws.Range("C217").Value = [double]'0.5256449444640776'
$ws.Range("D217").Value = [double]'0.4771063166343456'
$ws.Range("E217").Value = [double]'0.6332864169154406'
$ws.Range("F217").Value = [double]'0.4259939490975391'
$ws.Range("G217").Value = [double]'3.442081942971508'
$ws.Range("H217").Value = 'TZP'

$ws.Range("A218").Value = 'HospitalCHBH:Ward_ED_ICUOther'
$ws.Range("B218").Value = [double]'3.481924044116067'
$ws.Range("C218").Value = [double]'0.2413271077402592'
$ws.Range("D218").Value = [double]'5.169684578427574'
$ws.Range("E218").Value = [double]'2.344894330308487e-07'
$ws.Range("F218").Value = [double]'2.177119628438133'
$ws.Range("G218").Value = [double]'5.611068249551265'
$ws.Range("H218").Value = 'TZP'

$ws.Range("A219").Value = 'HospitalCNGMO:Ward_ED_ICUOther'
$ws.Range("B219").Value = [double]'3.964478548993912'
$ws.Range("C219").Value = [double]'0.497219382080847'
$ws.Range("D219").Value = [double]'2.77015414640372'
$ws.Range("E219").Value = [double]'0.005602976910376392'
$ws.Range("F219").Value = [double]'1.495851443367732'
$ws.Range("G219").Value = [double]'10.6768159769677'
$ws.Range("H219").Value = 'TZP'

$ws.Range("A220").Value = 'HospitalRabta:Ward_ED_ICUOther'
$ws.Range("B220").Value = [double]'1.342512881943767'
$ws.Range("C220").Value = [double]'0.1952084977272594'
$ws.Range("D220").Value = [double]'1.508864349006028'
$ws.Range("E220").Value = [double]'0.1313334521284583'
$ws.Range("F220").Value = [double]'0.9155959620652931'
$ws.Range("G220").Value = [double]'1.968887533491604'
$ws.Range("H220").Value = 'TZP'

$ws.Range("A221").Value = 'HospitalTCB:Ward_ED_ICUOther'
$ws.Range("B221").Value = [double]'3.440137586009126'
$ws.Range("C221").Value = [double]'0.1840665835476588'
$ws.Range("D221").Value = [double]'6.712307267866879'
$ws.Range("E221").Value = [double]'1.915707510172705e-11'
$ws.Range("F221").Value = [double]'2.399947359547067'
$ws.Range("G221").Value = [double]'4.940186492584962'
$ws.Range("H221").Value = 'TZP'

